$d = $word.ActiveDocument

$pairs = @(
    @("82×61=5002", "76×81=6156"),
    @("16×91=1456", "98×16=1568"),
    @("49×40=1960", "35×43=1505"),
    @("53×39=2067", "95×15=1425"),
    @("62×38=2356", "96×39=3744"),
    @("81×19=1539", "54×78=4212"),
    @("49×38=1862", "92×69=6348"),
    @("12×45=540",  "36×27=972"),
    @("81×58=4698", "14×24=336"),
    @("25×79=1975", "23×74=1702"),
    @("59×37=2183", "94×22=2068"),
    @("83×88=7304", "98×79=7742"),
    @("30×87=2610", "50×76=3800"),
    @("25×35=875",  "92×21=1932"),
    @("92×98=9016", "13×93=1209"),
    @("76×22=1672", "50×50=2500"),
    @("78×46=3588", "78×47=3666"),
    @("92×44=4048", "81×32=2592"),
    @("85×66=5610", "91×74=6734"),
    @("98×77=7546", "82×82=6724"),
    @("20×18=360",  "86×68=5848"),
    @("62×71=4402", "59×73=4307"),
    @("40×75=3000", "40×28=1120"),
    @("22×93=2046", "96×29=2784"),
    @("69×52=3588", "38×77=2926")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
